$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 2606.416448354451
$ws.Cells.Item(3, 2).Value = 3053.31058337847
$ws.Cells.Item(4, 2).Value = 2847.640197843274
$ws.Cells.Item(5, 2).Value = 2459.147967212766
$ws.Cells.Item(6, 2).Value = 2248.144178638431
$ws.Cells.Item(7, 2).Value = 2257.849154695007
$ws.Cells.Item(8, 2).Value = 2603.907684341763
$ws.Cells.Item(9, 2).Value = 2025.798166228645
$ws.Cells.Item(10, 2).Value = 2390.981695750565
$ws.Cells.Item(11, 2).Value = 2727.393374112728
$ws.Cells.Item(12, 2).Value = 2413.774641403596
$ws.Cells.Item(13, 2).Value = 2089.423576260539
$ws.Cells.Item(14, 2).Value = 2561.478901755176
$ws.Cells.Item(15, 2).Value = 2627.327569072482
$ws.Cells.Item(16, 2).Value = 2730.563470972452
$ws.Cells.Item(17, 2).Value = 2696.586379383648
$ws.Cells.Item(18, 2).Value = 2550.983605975811
$ws.Cells.Item(19, 2).Value = 2600.14071556693
$ws.Cells.Item(20, 2).Value = 2420.024288839112
$ws.Cells.Item(21, 2).Value = 2154.203252790418
$ws.Cells.Item(22, 2).Value = 2305.650715956946
$ws.Cells.Item(23, 2).Value = 2881.905875171302
$ws.Cells.Item(24, 2).Value = 3015.5767836972
$ws.Cells.Item(25, 2).Value = 2606.543628668703
$ws.Cells.Item(26, 2).Value = 2447.78085163456
$ws.Cells.Item(27, 2).Value = 2528.232624638157
$ws.Cells.Item(28, 2).Value = 2323.480365554033
$ws.Cells.Item(29, 2).Value = 2403.23353332432
$ws.Cells.Item(30, 2).Value = 2706.704806540611
$ws.Cells.Item(31, 2).Value = 2471.872769223498
$ws.Cells.Item(32, 2).Value = 2539.394945461609
$ws.Cells.Item(33, 2).Value = 3027.039581500004
$ws.Cells.Item(34, 2).Value = 2776.507774859863
$ws.Cells.Item(35, 2).Value = 2494.761332372781
$ws.Cells.Item(36, 2).Value = 2617.425363158868
$ws.Cells.Item(37, 2).Value = 2408.512815481193
$ws.Cells.Item(38, 2).Value = 2722.326811720019
$ws.Cells.Item(39, 2).Value = 2374.890503470433
$ws.Cells.Item(40, 2).Value = 1703.402102931034
$ws.Cells.Item(41, 2).Value = 2580.988056903566
$ws.Cells.Item(42, 2).Value = 2450.994004403067
$ws.Cells.Item(43, 2).Value = 3108.898085842262
$ws.Cells.Item(44, 2).Value = 2470.12824808735
$ws.Cells.Item(45, 2).Value = 2809.194562079102
$ws.Cells.Item(46, 2).Value = 2720.817157426522
$ws.Cells.Item(47, 2).Value = 2686.275312401783
$ws.Cells.Item(48, 2).Value = 2436.191105308698
$ws.Cells.Item(49, 2).Value = 2227.762205897822
$ws.Cells.Item(50, 2).Value = 2132.204060600924
$ws.Cells.Item(51, 2).Value = 2386.987192811604
$ws.Cells.Item(52, 2).Value = 2234.665633016381
$ws.Cells.Item(53, 2).Value = 2736.717304436537
$ws.Cells.Item(54, 2).Value = 2503.349912097522
$ws.Cells.Item(55, 2).Value = 2659.145710987831
$ws.Cells.Item(56, 2).Value = 2833.446202761741
$ws.Cells.Item(57, 2).Value = 2568.582053379174
$ws.Cells.Item(58, 2).Value = 2605.053143552518
$ws.Cells.Item(59, 2).Value = 2466.740220344659
$ws.Cells.Item(60, 2).Value = 2453.699961771757
$ws.Cells.Item(61, 2).Value = 2619.109955049199
$ws.Cells.Item(62, 2).Value = 2087.589981833594
$ws.Cells.Item(63, 2).Value = 2465.773663326734
$ws.Cells.Item(64, 2).Value = 2420.125491387243
$ws.Cells.Item(65, 2).Value = 2871.126459954377
$ws.Cells.Item(66, 2).Value = 2488.411217610107
$ws.Cells.Item(67, 2).Value = 2702.28057508017
$ws.Cells.Item(68, 2).Value = 2834.165304847314
$ws.Cells.Item(69, 2).Value = 2874.225766977246
$ws.Cells.Item(70, 2).Value = 2329.805557297021
$ws.Cells.Item(71, 2).Value = 2596.899551695898
$ws.Cells.Item(72, 2).Value = 2535.842863725175
$ws.Cells.Item(73, 2).Value = 2476.109109150387
$ws.Cells.Item(74, 2).Value = 2635.381609629572
$ws.Cells.Item(75, 2).Value = 2666.3581459758
$ws.Cells.Item(76, 2).Value = 2464.991940481308
$ws.Cells.Item(77, 2).Value = 2667.515458255037
$ws.Cells.Item(78, 2).Value = 2261.297251888212
$ws.Cells.Item(79, 2).Value = 2156.732763888876
$ws.Cells.Item(80, 2).Value = 2380.855036120035
$ws.Cells.Item(81, 2).Value = 2204.476163855022
$ws.Cells.Item(82, 2).Value = 2419.577064306668
$ws.Cells.Item(83, 2).Value = 2530.071181170762
$ws.Cells.Item(84, 2).Value = 2676.690950262805
$ws.Cells.Item(85, 2).Value = 2686.79899209044
$ws.Cells.Item(86, 2).Value = 2515.713103419027
$ws.Cells.Item(87, 2).Value = 2173.736883585239
$ws.Cells.Item(88, 2).Value = 2463.656135190604
$ws.Cells.Item(89, 2).Value = 2417.50804182963
$ws.Cells.Item(90, 2).Value = 2431.76083673618
$ws.Cells.Item(91, 2).Value = 2126.855395295528
$ws.Cells.Item(92, 2).Value = 2442.016240284512
$ws.Cells.Item(93, 2).Value = 2073.189548946541
$ws.Cells.Item(94, 2).Value = 2808.052179550587
$ws.Cells.Item(95, 2).Value = 2599.775441453516
$ws.Cells.Item(96, 2).Value = 2682.618647144555
$ws.Cells.Item(97, 2).Value = 2247.319420936009
$ws.Cells.Item(98, 2).Value = 2035.384321698455
$ws.Cells.Item(99, 2).Value = 2779.233401436265
$ws.Cells.Item(100, 2).Value = 2320.14829650485
$ws.Cells.Item(101, 2).Value = 2259.48632838731
